{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n// Applies the offer-letter text updates described by the diff:\n//   - offer date 23-09-2025 -> 24-09-2025\n//   - candidate name \"Virat Kohli\" -> \"Virat Kohli A\" (4 occurrences)\n//   - address line \"Surya Nagar,\" -> \"Sudama Nagar,\"\n//   - address/pin \"Maha Laxmi Nagar, Maharashtra, 452066\" -> \"...451111\"\n//   - phone number 9988776655 -> 9879809878\n//   - job title \"Solution Architect\" -> \"Team Lead\" (2 occurrences)\n//   - joining date 30-08-2025 -> 25-09-2025\n\nasync function replaceAll(context, findText, replaceText) {\n  const results = context.document.body.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nawait replaceAll(context, \"23-09-2025\", \"24-09-2025\");\nawait replaceAll(context, \"Virat Kohli\", \"Virat Kohli A\");\nawait replaceAll(context, \"Surya Nagar,\", \"Sudama Nagar,\");\nawait replaceAll(\n  context,\n  \"Maha Laxmi Nagar, Maharashtra, 452066\",\n  \"Maha Laxmi Nagar, Maharashtra, 451111\"\n);\nawait replaceAll(context, \"9988776655\", \"9879809878\");\nawait replaceAll(context, \"Solution Architect\", \"Team Lead\");\nawait replaceAll(context, \"30-08-2025\", \"25-09-2025\");\n", "ps1": "# Word COM interop script: apply the offer-letter text updates.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute(\n        $findText,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n}\n\n# Offer date at top of the letter.\nReplace-AllText \"23-09-2025\" \"24-09-2025\"\n\n# Candidate name -> \"Virat Kohli A\" (appears 4 times: address block, greeting,\n# signature block, and summary table).\nReplace-AllText \"Virat Kohli\" \"Virat Kohli A\"\n\n# Candidate address lines.\nReplace-AllText \"Surya Nagar,\" \"Sudama Nagar,\"\nReplace-AllText \"Maha Laxmi Nagar, Maharashtra, 452066\" \"Maha Laxmi Nagar, Maharashtra, 451111\"\n\n# Candidate phone number.\nReplace-AllText \"9988776655\" \"9879809878\"\n\n# Job title -> \"Team Lead\" (appears in the offer paragraph and the summary\n# table). After the replace-all, the engine can fold the trailing space run\n# that follows \"Solution Architect\" in the offer paragraph into the newly\n# replaced run (both share the same bold run formatting). Re-toggle Bold on\n# that trailing space to force Word to keep it as its own run again, just\n# like it was before the edit.\nReplace-AllText \"Solution Architect\" \"Team Lead\"\n\n$fixRange = $d.Content\n$fixRange.Find.ClearFormatting()\n$fixRange.Find.Execute(\"Team Lead \")\nif ($fixRange.Find.Found) {\n    $trailingSpace = $d.Range($fixRange.End - 1, $fixRange.End)\n    if ($trailingSpace.Text -eq \" \") {\n        $trailingSpace.Bold = 0\n        $trailingSpace.Bold = 1\n    }\n}\n\n# Joining date.\nReplace-AllText \"30-08-2025\" \"25-09-2025\"\n"}
